$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula-driven cells K8/L8 (row 8) and K9/L9 (row 9) with
# their plain computed constant values (formulas removed, results kept).
$ws.Range("K8").Value = 0.132
$ws.Range("L8").Value = 0.21296
$ws.Range("K9").Value = 0.132
$ws.Range("L9").Value = 0.21296

# Move the sheet's scroll position / active selection to F8.
$null = $ws.Range("F8").Select()
